# Applies the "Added the coulomb potential" edit to the workbook.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws4 = $wb.Worksheets.Item("Sheet4")

# Fix the mislabeled header on Sheet2 (H1): "COULOMBCOUPLING" -> "COULOMB COUPLING"
$ws2.Range("H1").Value = "COULOMB COUPLING"

# Sheet3: halve the isospin (F) coupling values for the coulomb potential
$ws3.Range("F2").Value = 0.5
$ws3.Range("F3").Value = 0.5
$ws3.Range("F4").Value = 0.5

# Sheet4: halve the (negative) isospin (F) coupling values for the coulomb potential
$ws4.Range("F2").Value = -0.5
$ws4.Range("F3").Value = -0.5
$ws4.Range("F4").Value = -0.5

# Update each sheet's selection/active-cell state
$ws1.Activate()
$ws1.Range("D17").Select()

$ws2.Activate()
$ws2.Range("F9").Select()

$ws3.Activate()
$ws3.Range("F2:F4").Select()

# Sheet4 ends up the active tab in the saved workbook
$ws4.Activate()
$ws4.Range("E21").Select()
